$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 5838.4062
$ws.Range("I15").Value = 5838.4062
$ws.Range("K15").Value = 17515.2186
$ws.Range("M15").Value = -17346.2186
$ws.Range("H32").Value = 25003046
$ws.Range("I32").Value = 40003936
$ws.Range("J32").Value = 14288126
$ws.Range("K32").Value = 40003936
$ws.Range("L32").Value = 14288126
$ws.Range("M32").Value = -40003610
$ws.Range("N32").Value = -14288778
$ws.Range("H51").Value = 4087.7778
$ws.Range("J51").Value = 4184.2
$ws.Range("L51").Value = 4184.2
$ws.Range("N51").Value = -5152.2
$ws.Range("H64").Value = 4245.4
$ws.Range("J64").Value = 5182.857
$ws.Range("L64").Value = 5182.857
$ws.Range("N64").Value = -5678.857
$ws.Range("H67").Value = 4245.4
$ws.Range("J67").Value = 5182.857
$ws.Range("L67").Value = 5182.857
$ws.Range("N67").Value = -6898.857
$ws.Range("H86").Value = 117649430
$ws.Range("I86").Value = 200002220
$ws.Range("J86").Value = 2601.2856
$ws.Range("K86").Value = 200002220
$ws.Range("L86").Value = 2601.2856
$ws.Range("M86").Value = -200001097
$ws.Range("N86").Value = -4847.2856
$ws.Range("H89").Value = 117649430
$ws.Range("I89").Value = 200002220
$ws.Range("J89").Value = 2601.2856
$ws.Range("K89").Value = 1000011100
$ws.Range("L89").Value = 13006.428
$ws.Range("M89").Value = -1000005484
$ws.Range("N89").Value = -24238.428
$ws.Range("H98").Value = 592
$ws.Range("I98").Value = 421.6
$ws.Range("K98").Value = 421.6
$ws.Range("M98").Value = 1076.4
$ws.Range("H106").Value = 25723840
$ws.Range("I106").Value = 33436402
$ws.Range("K106").Value = 33436402
$ws.Range("M106").Value = -33435771
$ws.Range("H115").Value = 62504140
$ws.Range("I115").Value = 66670748
$ws.Range("K115").Value = 200012244
$ws.Range("M115").Value = -200010677
$ws.Range("H122").Value = 592
$ws.Range("I122").Value = 421.6
$ws.Range("K122").Value = 1264.8
$ws.Range("M122").Value = 1185.2
$ws.Range("H125").Value = 1941.5454
$ws.Range("J125").Value = 4993.6665
$ws.Range("L125").Value = 44942.9985
$ws.Range("N125").Value = -49862.9985
$ws.Range("H131").Value = 3490.5386
$ws.Range("I131").Value = 1604.2106
$ws.Range("J131").Value = 8610.571
$ws.Range("K131").Value = 4812.6318
$ws.Range("L131").Value = 25831.713
$ws.Range("M131").Value = 227.3681999999999
$ws.Range("N131").Value = -35911.713
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("N133").ClearContents()
$ws.Range("H135").Value = 1079.1082
$ws.Range("I135").Value = 546.5517
$ws.Range("K135").Value = 4918.9653
$ws.Range("M135").Value = -2383.9653
$ws.Range("H137").Value = 3339.082
$ws.Range("I137").Value = 3934.1428
$ws.Range("J137").Value = 2023.6842
$ws.Range("K137").Value = 11802.4284
$ws.Range("L137").Value = 6071.0526
$ws.Range("M137").Value = -9252.428400000001
$ws.Range("N137").Value = -11171.0526
$ws.Range("H138").Value = 10640150
$ws.Range("I138").Value = 31251448
$ws.Range("K138").Value = 93754344
$ws.Range("M138").Value = -93749204
$ws.Range("H141").Value = 5996.303
$ws.Range("I141").Value = 5714.174
$ws.Range("K141").Value = 17142.522
$ws.Range("M141").Value = -11962.522

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1654.3334
$ws.Range("I4").Value = 963
$ws.Range("K4").Value = 963
$ws.Range("M4").Value = -847
$ws.Range("H11").Value = 6428.75
$ws.Range("I11").Value = 5524.5
$ws.Range("J11").Value = 7333
$ws.Range("K11").Value = 5524.5
$ws.Range("L11").Value = 7333
$ws.Range("M11").Value = -5380.5
$ws.Range("N11").Value = -7621
$ws.Range("H12").Value = 1873.2
$ws.Range("I12").Value = 675
$ws.Range("J12").Value = 6666
$ws.Range("K12").Value = 675
$ws.Range("L12").Value = 6666
$ws.Range("M12").Value = -502
$ws.Range("N12").Value = -7012
$ws.Range("H32").Value = 5070.0894
$ws.Range("I32").Value = 5103.442
$ws.Range("J32").Value = 4959.769
$ws.Range("K32").Value = 5103.442
$ws.Range("L32").Value = 4959.769
$ws.Range("M32").Value = -4816.442
$ws.Range("N32").Value = -5533.769
$ws.Range("H45").Value = 6231.091
$ws.Range("I45").Value = 8517.368
$ws.Range("J45").Value = 3128.2856
$ws.Range("K45").Value = 8517.368
$ws.Range("L45").Value = 3128.2856
$ws.Range("M45").Value = -8140.368
$ws.Range("N45").Value = -3882.2856
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H54").Value = 400000
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H58").Value = 43333
$ws.Range("J58").Value = 43333
$ws.Range("L58").Value = 43333
$ws.Range("N58").Value = -44193
$ws.Range("H61").Value = 4833.1816
$ws.Range("I61").Value = 5077.6553
$ws.Range("J61").Value = 3060.75
$ws.Range("K61").Value = 5077.6553
$ws.Range("L61").Value = 3060.75
$ws.Range("M61").Value = -4865.6553
$ws.Range("N61").Value = -3484.75
$ws.Range("H74").Value = 6733.606
$ws.Range("I74").Value = 6660.271
$ws.Range("K74").Value = 6660.271
$ws.Range("M74").Value = -5786.271
$ws.Range("H77").Value = 6733.606
$ws.Range("I77").Value = 6660.271
$ws.Range("K77").Value = 33301.355
$ws.Range("M77").Value = -28933.355
$ws.Range("H102").Value = 2986.7144
$ws.Range("I102").Value = 2790.739
$ws.Range("K102").Value = 2790.739
$ws.Range("M102").Value = -1168.739
$ws.Range("H122").Value = 4914.2
$ws.Range("I122").Value = 2012
$ws.Range("K122").Value = 6036
$ws.Range("M122").Value = -3586
$ws.Range("H132").Value = 1151.4546
$ws.Range("I132").Value = 1016.6
$ws.Range("K132").Value = 3049.8
$ws.Range("M132").Value = -519.8000000000002
$ws.Range("H136").Value = 4833.1816
$ws.Range("I136").Value = 5077.6553
$ws.Range("J136").Value = 3060.75
$ws.Range("K136").Value = 15232.9659
$ws.Range("L136").Value = 9182.25
$ws.Range("M136").Value = -12682.9659
$ws.Range("N136").Value = -14282.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 6500
$ws.Range("J15").Value = 6500
$ws.Range("L15").Value = 6500
$ws.Range("N15").Value = -6954
$ws.Range("H19").Value = 1150
$ws.Range("J19").Value = 2000
$ws.Range("L19").Value = 2000
$ws.Range("N19").Value = -2346
$ws.Range("H20").Value = 2401.8125
$ws.Range("J20").Value = 2507
$ws.Range("L20").Value = 2507
$ws.Range("N20").Value = -3001
$ws.Range("H86").Value = 2513.5789
$ws.Range("I86").Value = 2340.0833
$ws.Range("J86").Value = 2811
$ws.Range("K86").Value = 2340.0833
$ws.Range("L86").Value = 2811
$ws.Range("M86").Value = -1217.0833
$ws.Range("N86").Value = -5057
$ws.Range("H89").Value = 2513.5789
$ws.Range("I89").Value = 2340.0833
$ws.Range("J89").Value = 2811
$ws.Range("K89").Value = 11700.4165
$ws.Range("L89").Value = 14055
$ws.Range("M89").Value = -6084.416499999999
$ws.Range("N89").Value = -25287
$ws.Range("H134").Value = 1936.541
$ws.Range("I134").Value = 1953.8334
$ws.Range("J134").Value = 899
$ws.Range("K134").Value = 5861.5002
$ws.Range("L134").Value = 2697
$ws.Range("M134").Value = -3326.5002
$ws.Range("N134").Value = -7767

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 160773.5
$ws.Range("J20").Value = 160773.5
$ws.Range("L20").Value = 160773.5
$ws.Range("N20").Value = -161245.5
$ws.Range("H30").Value = 160773.5
$ws.Range("J30").Value = 160773.5
$ws.Range("L30").Value = 160773.5
$ws.Range("N30").Value = -160955.5
$ws.Range("H31").Value = 1883.8636
$ws.Range("I31").Value = 1188.5625
$ws.Range("J31").Value = 3738
$ws.Range("K31").Value = 1188.5625
$ws.Range("L31").Value = 3738
$ws.Range("M31").Value = -893.5625
$ws.Range("N31").Value = -4328
$ws.Range("H34").Value = 1883.8636
$ws.Range("I34").Value = 1188.5625
$ws.Range("J34").Value = 3738
$ws.Range("K34").Value = 1188.5625
$ws.Range("L34").Value = 3738
$ws.Range("M34").Value = -986.5625
$ws.Range("N34").Value = -4142
$ws.Range("H94").Value = 5444.625
$ws.Range("J94").Value = 5081
$ws.Range("L94").Value = 5081
$ws.Range("N94").Value = -5983
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H99").Value = 13588.486
$ws.Range("I99").Value = 10422.714
$ws.Range("J99").Value = 15515.479
$ws.Range("K99").Value = 10422.714
$ws.Range("L99").Value = 15515.479
$ws.Range("M99").Value = -8924.714
$ws.Range("N99").Value = -18511.479
$ws.Range("H107").Value = 1480.2222
$ws.Range("I107").Value = 1036.4286
$ws.Range("J107").Value = 3033.5
$ws.Range("K107").Value = 1036.4286
$ws.Range("L107").Value = 3033.5
$ws.Range("M107").Value = 883.5714
$ws.Range("N107").Value = -6873.5
$ws.Range("H122").Value = 6269.3335
$ws.Range("I122").Value = 6237.923
$ws.Range("K122").Value = 18713.769
$ws.Range("M122").Value = -16263.769
$ws.Range("H126").Value = 13588.486
$ws.Range("I126").Value = 10422.714
$ws.Range("J126").Value = 15515.479
$ws.Range("K126").Value = 31268.142
$ws.Range("L126").Value = 46546.437
$ws.Range("M126").Value = -28798.142
$ws.Range("N126").Value = -51486.437
$ws.Range("H128").Value = 160773.5
$ws.Range("J128").Value = 160773.5
$ws.Range("L128").Value = 160773.5
$ws.Range("N128").Value = -170733.5
$ws.Range("H129").Value = 57714.07
$ws.Range("J129").Value = 57714.07
$ws.Range("L129").Value = 57714.07
$ws.Range("N129").Value = -67714.07000000001
$ws.Range("H130").Value = 200000
$ws.Range("J130").Value = 200000
$ws.Range("L130").Value = 200000
$ws.Range("N130").Value = -210040
$ws.Range("H132").Value = 10628.053
$ws.Range("I132").Value = 3612.1538
$ws.Range("J132").Value = 25829.166
$ws.Range("K132").Value = 10836.4614
$ws.Range("L132").Value = 77487.49800000001
$ws.Range("M132").Value = -8306.4614
$ws.Range("N132").Value = -82547.49800000001
$ws.Range("H134").Value = 3254.121
$ws.Range("I134").Value = 3506.963
$ws.Range("J134").Value = 2116.3333
$ws.Range("K134").Value = 10520.889
$ws.Range("L134").Value = 6348.999899999999
$ws.Range("M134").Value = -7985.889000000001
$ws.Range("N134").Value = -11418.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 77.7619
$ws.Range("I2").Value = 21.181818
$ws.Range("J2").Value = 140
$ws.Range("K2").Value = 127.090908
$ws.Range("L2").Value = 840
$ws.Range("M2").Value = -14.090908
$ws.Range("N2").Value = -1066
$ws.Range("H11").Value = 3578.1
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 3578.1
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 10734.3
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -11014.3
$ws.Range("H14").Value = 392
$ws.Range("I14").Value = 392
$ws.Range("K14").Value = 1176
$ws.Range("M14").Value = -1003
$ws.Range("H17").Value = 381.85715
$ws.Range("I17").Value = 447.5
$ws.Range("K17").Value = 1342.5
$ws.Range("M17").Value = -1173.5
$ws.Range("H26").Value = 1104.75
$ws.Range("I26").Value = 172.33333
$ws.Range("J26").Value = 3902
$ws.Range("K26").Value = 516.99999
$ws.Range("L26").Value = 11706
$ws.Range("M26").Value = -228.99999
$ws.Range("N26").Value = -12282
$ws.Range("H38").Value = 427.6875
$ws.Range("J38").Value = 556.5
$ws.Range("L38").Value = 1669.5
$ws.Range("N38").Value = -2363.5
$ws.Range("H39").Value = 6503.8184
$ws.Range("J39").Value = 9248.857
$ws.Range("L39").Value = 27746.571
$ws.Range("N39").Value = -28334.571
$ws.Range("H103").Value = 299.875
$ws.Range("I103").Value = 199.85715
$ws.Range("J103").Value = 1000
$ws.Range("K103").Value = 599.5714499999999
$ws.Range("L103").Value = 3000
$ws.Range("M103").Value = 279.4285500000001
$ws.Range("N103").Value = -4758
$ws.Range("H105").Value = 14929
$ws.Range("J105").Value = 14929
$ws.Range("L105").Value = 44787
$ws.Range("N105").Value = -50029
$ws.Range("H107").Value = 498.35715
$ws.Range("J107").Value = 536.75
$ws.Range("L107").Value = 1610.25
$ws.Range("N107").Value = -5450.25
$ws.Range("H113").Value = 2264.2222
$ws.Range("J113").Value = 2234.75
$ws.Range("L113").Value = 6704.25
$ws.Range("N113").Value = -11044.25
$ws.Range("H116").Value = 5668135
$ws.Range("I116").Value = 5668135
$ws.Range("K116").Value = 17004405
$ws.Range("M116").Value = -17000963
$ws.Range("H117").Value = 1585.4445
$ws.Range("I117").Value = 571.75
$ws.Range("J117").Value = 2396.4
$ws.Range("K117").Value = 1715.25
$ws.Range("L117").Value = 7189.200000000001
$ws.Range("M117").Value = 1726.75
$ws.Range("N117").Value = -14073.2
$ws.Range("H132").Value = 4949.3335
$ws.Range("I132").Value = 4804
$ws.Range("J132").Value = 4978.4
$ws.Range("K132").Value = 43236
$ws.Range("L132").Value = 44805.6
$ws.Range("M132").Value = -40706
$ws.Range("N132").Value = -49865.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3701.5757
$ws.Range("I102").Value = 3962.682
$ws.Range("K102").Value = 3962.682
$ws.Range("M102").Value = -2340.682
$ws.Range("H113").Value = 5924.25
$ws.Range("I113").Value = 5850
$ws.Range("K113").Value = 5850
$ws.Range("M113").Value = -3680
$ws.Range("H126").Value = 5969.8237
$ws.Range("I126").Value = 5624.3335
$ws.Range("K126").Value = 16873.0005
$ws.Range("M126").Value = -14403.0005
$ws.Range("H132").Value = 3547.375
$ws.Range("I132").Value = 3547.375
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10642.125
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8112.125
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2744.1667
$ws.Range("I7").Value = 2726.3333
$ws.Range("J7").Value = 2797.6667
$ws.Range("K7").Value = 2726.3333
$ws.Range("L7").Value = 2797.6667
$ws.Range("M7").Value = -2614.3333
$ws.Range("N7").Value = -3021.6667
$ws.Range("H16").Value = 2843224.8
$ws.Range("I16").Value = 3573389.2
$ws.Range("J16").Value = 3696.889
$ws.Range("K16").Value = 3573389.2
$ws.Range("L16").Value = 3696.889
$ws.Range("M16").Value = -3573219.2
$ws.Range("N16").Value = -4036.889
$ws.Range("H40").Value = 6887.15
$ws.Range("I40").Value = 6889.769
$ws.Range("K40").Value = 6889.769
$ws.Range("M40").Value = -6753.769
$ws.Range("H46").Value = 3437.1
$ws.Range("I46").Value = 2500
$ws.Range("J46").Value = 4061.8333
$ws.Range("K46").Value = 2500
$ws.Range("L46").Value = 4061.8333
$ws.Range("M46").Value = -2312
$ws.Range("N46").Value = -4437.8333
$ws.Range("H68").Value = 8550680
$ws.Range("I68").Value = 9805074
$ws.Range("J68").Value = 20799
$ws.Range("K68").Value = 9805074
$ws.Range("L68").Value = 20799
$ws.Range("M68").Value = -9804325
$ws.Range("N68").Value = -22297
$ws.Range("H71").Value = 8550680
$ws.Range("I71").Value = 9805074
$ws.Range("J71").Value = 20799
$ws.Range("K71").Value = 49025370
$ws.Range("L71").Value = 103995
$ws.Range("M71").Value = -49021626
$ws.Range("N71").Value = -111483
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H122").Value = 3196.1538
$ws.Range("I122").Value = 3390.05
$ws.Range("K122").Value = 10170.15
$ws.Range("M122").Value = -7720.150000000001
$ws.Range("H126").Value = 2744.1667
$ws.Range("I126").Value = 2726.3333
$ws.Range("J126").Value = 2797.6667
$ws.Range("K126").Value = 8178.999899999999
$ws.Range("L126").Value = 8393.000100000001
$ws.Range("M126").Value = -5708.999899999999
$ws.Range("N126").Value = -13333.0001
$ws.Range("H132").Value = 1000004
$ws.Range("I132").Value = 1000004
$ws.Range("K132").Value = 3000012
$ws.Range("M132").Value = -2997482
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H135").Value = 59998
$ws.Range("J135").Value = 59998
$ws.Range("L135").Value = 59998
$ws.Range("N135").Value = -70138
$ws.Range("H136").Value = 100001.14
$ws.Range("I136").Value = 100001.5
$ws.Range("J136").Value = 99999
$ws.Range("K136").Value = 300004.5
$ws.Range("L136").Value = 299997
$ws.Range("M136").Value = -297454.5
$ws.Range("N136").Value = -305097

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2651.6
$ws.Range("I122").Value = 2769.5715
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 8308.7145
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -5858.7145
$ws.Range("N122").Value = -7900
$ws.Range("H126").Value = 7368.55
$ws.Range("J126").Value = 14200.6
$ws.Range("L126").Value = 42601.8
$ws.Range("N126").Value = -47541.8
$ws.Range("H132").Value = 1975.5957
$ws.Range("I132").Value = 1953.326
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 5859.978
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -3329.978
$ws.Range("N132").Value = -14060
